# Automatic update of files.
#
# The underlying dataset rows (20-22, 28-29, 32-35) got re-sorted: the
# observation records (same columns, same sheet) moved between row numbers.
# This script rewrites each affected row's cells in place so the sheet ends
# up holding the same per-row field values the re-sort produced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 20 -----------------------------------------------------------
$ws.Range("A20").Value = 131092554
$ws.Range("B20").Value = 57884
$ws.Range("E20").Value = 100109
$ws.Range("F20").Value = "Tretåig hackspett"
$ws.Range("G20").Value = "Picoides tridactylus"
$ws.Range("H20").Value = "(Linnaeus, 1758)"
$ws.Range("M20").Value = "äldre spår"
$ws.Range("Q20").Value = 585147
$ws.Range("R20").Value = 7060312
$ws.Range("S20").Value = 15
$ws.Range("Z20").ClearContents()
$ws.Range("AB20").ClearContents()
$ws.Range("AC20").Value = "Äldre ringhack, gran"
$ws.Range("AW20").Value = "Daniel Rutschman"
$ws.Range("AX20").Value = "Daniel Rutschman"

# --- Row 21 -----------------------------------------------------------
$ws.Range("A21").Value = 131086957
$ws.Range("Q21").Value = 585162
$ws.Range("R21").Value = 7060573
$ws.Range("S21").Value = 10
$ws.Range("M21").Value = "färska spår"
$ws.Range("Z21").Value = "12:21"
$ws.Range("AB21").Value = "12:21"
$ws.Range("AC21").Value = "Ringhack på gran"
$ws.Range("AW21").Value = "Kim Hultgren"
$ws.Range("AX21").Value = "Kim Hultgren"

# --- Row 22 -----------------------------------------------------------
$ws.Range("A22").Value = 131092560
$ws.Range("B22").Value = 91804
$ws.Range("E22").Value = 1108
$ws.Range("F22").Value = "Harticka"
$ws.Range("G22").Value = "Pelloporus leporinus"
$ws.Range("H22").Value = "(Fr.) Krieglst."
$ws.Range("M22").ClearContents()
$ws.Range("Q22").Value = 585129
$ws.Range("R22").Value = 7060254
$ws.Range("Z22").Value = "15:17"
$ws.Range("AB22").Value = "15:17"
$ws.Range("AC22").ClearContents()

# --- Row 28 -----------------------------------------------------------
$ws.Range("A28").Value = 131085171
$ws.Range("B28").Value = 91804
$ws.Range("E28").Value = 1108
$ws.Range("F28").Value = "Harticka"
$ws.Range("G28").Value = "Pelloporus leporinus"
$ws.Range("H28").Value = "(Fr.) Krieglst."
$ws.Range("Q28").Value = 585222
$ws.Range("R28").Value = 7060254
$ws.Range("S28").Value = 15
$ws.Range("Z28").ClearContents()
$ws.Range("AB28").ClearContents()
$ws.Range("AW28").Value = "Daniel Rutschman"
$ws.Range("AX28").Value = "Daniel Rutschman"

# --- Row 29 -----------------------------------------------------------
$ws.Range("A29").Value = 131085178
$ws.Range("B29").Value = 91828
$ws.Range("E29").Value = 5432
$ws.Range("F29").Value = "Granticka"
$ws.Range("G29").Value = "Porodaedalea chrysoloma s.lat."
$ws.Range("H29").ClearContents()
$ws.Range("Q29").Value = 585225
$ws.Range("R29").Value = 7060258
$ws.Range("S29").Value = 10
$ws.Range("Z29").Value = "11:08"
$ws.Range("AB29").Value = "11:08"
$ws.Range("AW29").Value = "Kim Hultgren"
$ws.Range("AX29").Value = "Kim Hultgren"

# --- Row 32 -----------------------------------------------------------
$ws.Range("A32").Value = 131092585
$ws.Range("B32").Value = 91804
$ws.Range("E32").Value = 1108
$ws.Range("F32").Value = "Harticka"
$ws.Range("G32").Value = "Pelloporus leporinus"
$ws.Range("H32").Value = "(Fr.) Krieglst."
$ws.Range("Q32").Value = 585130
$ws.Range("R32").Value = 7060263

# --- Row 33 -----------------------------------------------------------
$ws.Range("A33").Value = 131085569
$ws.Range("Q33").Value = 585249
$ws.Range("R33").Value = 7060505

# --- Row 34 -----------------------------------------------------------
$ws.Range("A34").Value = 131087388
$ws.Range("Q34").Value = 585131
$ws.Range("R34").Value = 7060627
$ws.Range("S34").Value = 15
$ws.Range("Z34").ClearContents()
$ws.Range("AB34").ClearContents()
$ws.Range("AW34").Value = "Daniel Rutschman"
$ws.Range("AX34").Value = "Daniel Rutschman"

# --- Row 35 -----------------------------------------------------------
$ws.Range("A35").Value = 131092590
$ws.Range("B35").Value = 79243
$ws.Range("E35").Value = 6425
$ws.Range("F35").Value = "Garnlav"
$ws.Range("G35").Value = "Alectoria sarmentosa"
$ws.Range("H35").Value = "(Ach.) Ach."
$ws.Range("Q35").Value = 585145
$ws.Range("R35").Value = 7060230
$ws.Range("S35").Value = 10
$ws.Range("Z35").Value = "15:20"
$ws.Range("AB35").Value = "15:20"
$ws.Range("AW35").Value = "Kim Hultgren"
$ws.Range("AX35").Value = "Kim Hultgren"
